$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AH2").Value = -0.18200911947172416
$ws.Range("AI2").Value = 0.3634622042190456
$ws.Range("AN2").Value = -0.059974646419790154
$ws.Range("AO2").Value = 0.03358165110387526
$ws.Range("AP2").Value = -1.7859350105888387
$ws.Range("AH3").Value = -0.179855521060702
$ws.Range("AI3").Value = 0.31798466794609304
$ws.Range("AN3").Value = 0.02718499287085345
$ws.Range("AO3").Value = 0.033579680465720296
$ws.Range("AP3").Value = 0.8095667526856057
$ws.Range("AH4").Value = -0.17847722047095801
$ws.Range("AI4").Value = 0.27079195358971087
$ws.Range("AN4").Value = 0.11418113589614041
$ws.Range("AO4").Value = 0.03358816054668617
$ws.Range("AP4").Value = 3.3994459368334065
$ws.Range("AH5").Value = -0.1778465772415593
$ws.Range("AI5").Value = 0.2233215220040591
$ws.Range("AM5").Value = 4.849596732467987
$ws.Range("AN5").Value = 0.20101961152855202
$ws.Range("AO5").Value = 0.033607034117061105
$ws.Range("AP5").Value = 5.981474319583038
$ws.Range("AH6").Value = -0.17793552058426151
$ws.Range("AI6").Value = 0.17377794976946195
$ws.Range("AM6").Value = 3.7697921332936257
$ws.Range("AN6").Value = 0.28770633938880386
$ws.Range("AO6").Value = 0.033636246880734154
$ws.Range("AP6").Value = 8.553461401592742
$ws.Range("AH7").Value = -0.17871549829149863
$ws.Range("AI7").Value = 0.12299950289072481
$ws.Range("AM7").Value = 2.6630733803753195
$ws.Range("AN7").Value = 0.374247340620225
$ws.Range("AO7").Value = 0.03367574752246818
$ws.Range("AP7").Value = 11.113260080434156
$ws.Range("AH8").Value = -0.18015758990058472
$ws.Range("AI8").Value = 0.07098827562335741
$ws.Range("AM8").Value = 1.5309282182059158
$ws.Range("AN8").Value = 0.4606487140244784
$ws.Range("AO8").Value = 0.03372548774107696
$ws.Range("AP8").Value = 13.658771003137119
$ws.Range("AH9").Value = -0.18223246754213301
$ws.Range("AI9").Value = 0.01777678569098337
$ws.Range("AM9").Value = 0.38337292467449346
$ws.Range("AN9").Value = 0.5469166443179194
$ws.Range("AO9").Value = 0.03378542228999648
$ws.Range("AP9").Value = 16.187947559852045
$ws.Range("AH10").Value = -0.18491042708891078
$ws.Range("AI10").Value = -0.03663716481747331
$ws.Range("AM10").Value = -0.7901145500662903
$ws.Range("AN10").Value = 0.6330573955628669
$ws.Range("AO10").Value = 0.033855509011093646
$ws.Range("AP10").Value = 18.69880010829047
$ws.Range("AH11").Value = -0.18816140453281774
$ws.Range("AI11").Value = -0.09223746267211422
$ws.Range("AM11").Value = -1.9926131973518912
$ws.Range("AN11").Value = 0.7190773077139954
$ws.Range("AO11").Value = 0.03393570886718653
$ws.Range("AP11").Value = 21.189399948244283
$ws.Range("AH12").Value = -0.19189104450505803
$ws.Range("AI12").Value = -0.1489796626685081
$ws.Range("AM12").Value = -3.2293122434612735
$ws.Range("AN12").Value = 0.8047607487918693
$ws.Range("AO12").Value = 0.034025739509808
$ws.Range("AP12").Value = 23.651528530625914
$ws.Range("AH13").Value = -0.19573971972298657
$ws.Range("AI13").Value = -0.20672790880034916
$ws.Range("AM13").Value = -4.487938114687211
$ws.Range("AN13").Value = 0.888998465633821
$ws.Range("AO13").Value = 0.03412412088596913
$ws.Range("AP13").Value = 26.05190822657506
$ws.Range("AH14").Value = -0.19926330727813674
$ws.Range("AI14").Value = -0.2651986444535372
$ws.Range("AM14").Value = -5.788549644696605
$ws.Range("AN14").Value = 0.9704591510733902
$ws.Range("AO14").Value = 0.03422856620452665
$ws.Range("AP14").Value = 28.352316754215934
$ws.Range("AH15").Value = -0.20199753060877512
$ws.Range("AI15").Value = -0.3244155903469875
$ws.Range("AM15").Value = -7.123609989640499
$ws.Range("AN15").Value = 1.0478149377509935
$ws.Range("AO15").Value = 0.03433621956389443
$ws.Range("AP15").Value = 30.51631632891824
$ws.Range("AH16").Value = -0.20345928357878315
$ws.Range("AI16").Value = -0.38388371619035144
$ws.Range("AM16").Value = -8.476441051703024
$ws.Range("AN16").Value = 1.119741118893512
$ws.Range("AO16").Value = 0.034443720432976506
$ws.Range("AP16").Value = 32.50929646442807
$ws.Range("AH17").Value = -0.2031479880941241
$ws.Range("AI17").Value = -0.44392268030400905
$ws.Range("AM17").Value = -9.955141244059213
$ws.Range("AN17").Value = 1.1849158620175726
$ws.Range("AO17").Value = 0.03454729193114822
$ws.Range("AP17").Value = 34.29837176178893
$ws.Range("AH18").Value = -0.2005469837879335
$ws.Range("AI18").Value = -0.5038174443596343
$ws.Range("AM18").Value = -11.430289931876676
$ws.Range("AN18").Value = 1.2420199158701368
$ws.Range("AO18").Value = 0.03464285274927803
$ws.Range("AP18").Value = 35.852125829793884
$ws.Range("AH19").Value = -0.19512494826270507
$ws.Range("AI19").Value = -0.5639055044648694
$ws.Range("AM19").Value = -13.309821960835919
$ws.Range("AN19").Value = 1.2897363109251543
$ws.Range("AO19").Value = 0.03472615258545925
$ws.Range("AP19").Value = 37.14020168952436
$ws.Range("AH20").Value = -0.1863373463386389
$ws.Range("AI20").Value = -0.6235640039118161
$ws.Range("AM20").Value = -15.73311297377797
$ws.Range("AN20").Value = 1.3267500537633379
$ws.Range("AO20").Value = 0.03479293100158517
$ws.Range("AP20").Value = 38.132747531471004
$ws.Range("AH21").Value = -0.1736279067199838
$ws.Range("AI21").Value = -0.6828490490428352
$ws.Range("AM21").Value = -23.942451968241606
$ws.Range("AN21").Value = 1.3517478156699865
$ws.Range("AO21").Value = 0.03483909964008326
$ws.Range("AP21").Value = 38.79973448322891
$ws.Range("AH22").Value = -0.15643012445619756
$ws.Range("AI22").Value = -0.7412596003369537
$ws.Range("AM22").Value = -25.0
$ws.Range("AN22").Value = 1.363417615793142
$ws.Range("AO22").Value = 0.034860947773960926
$ws.Range("AP22").Value = 39.11017063086089
